# This script applies the Renaissance page-rank / ZGC / heap-1G docx update:
# - Rows 1-3 (1-indexed) get their summary numbers swapped for new "0M" placeholders
#   (their old values 82.12 / 203.52 / 1138 are relocated down into the detail rows
#   43-45, collapsing those multi-run tab-separated rows into a single value).
# - Row 4 total-time value updates 3479 -> 4245.
# - Several stat rows (6, 7, 8, 11, 12) get refreshed numbers.
# - Rows 43-45 (previously long tab-separated per-iteration breakdowns) collapse
#   down to the single relocated summary value.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $rowIndex, $newText) {
    $cellRange = $table.Cell($rowIndex, 1).Range
    $cellRange.Text = $newText
}

Set-CellText $t 1 "0M"
Set-CellText $t 2 "0M"
Set-CellText $t 3 "0M"
Set-CellText $t 4 "4245"

Set-CellText $t 6 "0.41208"
Set-CellText $t 7 "0.08319"
Set-CellText $t 8 "0.02179"

Set-CellText $t 11 "0.29635"
Set-CellText $t 12 "203.51944"

Set-CellText $t 44 "82.12"
Set-CellText $t 45 "203.52"
Set-CellText $t 46 "1138"
